# Generate Report for Archive
#
# 1. Update the "Status" text from "Ready for handoff" to "In Translation"
#    on all three sheets (Overview, zh-cn, de-de).
# 2. Shrink the "Status" column width on all three sheets from
#    17.2159881591797 to 13.4101845877511.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update status text -----------------------------------------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value     = "In Translation"
$dede.Range("C2").Value     = "In Translation"

# --- Shrink status column widths ---------------------------------------
# NOTE: Excel's ColumnWidth property snaps to a pixel grid (1/6-character
# steps on this host), so this is the closest value reachable through the
# object model to the canonical target width of 13.4101845877511.
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth     = 12.5
$dede.Range("C1").ColumnWidth     = 12.5
